# Auto-generated Excel COM-interop script applying numeric updates
# to the Phoenix_Profits workbook per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 469.92307
$ws.Range("J6").Value = 1050.5
$ws.Range("L6").Value = 3151.5
$ws.Range("N6").Value = -3375.5

$ws.Range("H116").Value = 7446.1387
$ws.Range("I116").Value = 6673.2383
$ws.Range("K116").Value = 6673.2383
$ws.Range("M116").Value = -3231.2383

$ws.Range("H132").Value = 2325.8604
$ws.Range("I132").Value = 2236.2563
$ws.Range("K132").Value = 6708.7689
$ws.Range("M132").Value = -4178.7689

$ws.Range("H135").Value = 1186.2727
$ws.Range("I135").Value = 1083.3334
$ws.Range("J135").Value = 1649.5
$ws.Range("K135").Value = 9750.000599999999
$ws.Range("L135").Value = 14845.5
$ws.Range("M135").Value = -7215.000599999999
$ws.Range("N135").Value = -19915.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1765.6842
$ws.Range("I2").Value = 1536.125
$ws.Range("J2").Value = 2990
$ws.Range("K2").Value = 1536.125
$ws.Range("L2").Value = 2990
$ws.Range("M2").Value = -1423.125
$ws.Range("N2").Value = -3216

$ws.Range("H16").Value = 5627.4287
$ws.Range("I16").Value = 3098.4
$ws.Range("K16").Value = 3098.4
$ws.Range("M16").Value = -2811.4

$ws.Range("H61").Value = 4915.5264
$ws.Range("I61").Value = 4043.889
$ws.Range("K61").Value = 4043.889
$ws.Range("M61").Value = -3831.889

$ws.Range("H116").Value = 1765.6842
$ws.Range("I116").Value = 1536.125
$ws.Range("J116").Value = 2990
$ws.Range("K116").Value = 1536.125
$ws.Range("L116").Value = 2990
$ws.Range("M116").Value = 757.875
$ws.Range("N116").Value = -7578

$ws.Range("H132").Value = 3642.4
$ws.Range("I132").Value = 4030.7144
$ws.Range("J132").Value = 2736.3333
$ws.Range("K132").Value = 12092.1432
$ws.Range("L132").Value = 8208.999899999999
$ws.Range("M132").Value = -9562.143199999999
$ws.Range("N132").Value = -13268.9999

$ws.Range("H136").Value = 4915.5264
$ws.Range("I136").Value = 4043.889
$ws.Range("K136").Value = 12131.667
$ws.Range("M136").Value = -9581.667000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1765.6842
$ws.Range("I3").Value = 1536.125
$ws.Range("J3").Value = 2990
$ws.Range("K3").Value = 1536.125
$ws.Range("L3").Value = 2990
$ws.Range("M3").Value = -1422.125
$ws.Range("N3").Value = -3218

$ws.Range("H134").Value = 7570.1333
$ws.Range("I134").Value = 4178.778
$ws.Range("K134").Value = 12536.334
$ws.Range("M134").Value = -10001.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 532.3077
$ws.Range("I47").Value = 157.27272
$ws.Range("J47").Value = 2595
$ws.Range("K47").Value = 471.81816
$ws.Range("L47").Value = 7785
$ws.Range("M47").Value = -40.81815999999998
$ws.Range("N47").Value = -8647

$ws.Range("H86").Value = 1654.3334
$ws.Range("I86").Value = 383
$ws.Range("J86").Value = 3243.5
$ws.Range("K86").Value = 1149
$ws.Range("L86").Value = 9730.5
$ws.Range("M86").Value = 37
$ws.Range("N86").Value = -12102.5

$ws.Range("H89").Value = 1654.3334
$ws.Range("I89").Value = 383
$ws.Range("J89").Value = 3243.5
$ws.Range("K89").Value = 3447
$ws.Range("L89").Value = 29191.5
$ws.Range("M89").Value = 2481
$ws.Range("N89").Value = -41047.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 127000
$ws.Range("J95").Value = 127000
$ws.Range("L95").Value = 127000
$ws.Range("N95").Value = -132492

$ws.Range("H102").Value = 5649.9585
$ws.Range("I102").Value = 5460.579
$ws.Range("K102").Value = 5460.579
$ws.Range("M102").Value = -3838.579

$ws.Range("H132").Value = 1503127.9
$ws.Range("I132").Value = 1582029.4
$ws.Range("K132").Value = 4746088.199999999
$ws.Range("M132").Value = -4743558.199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3010.75
$ws.Range("J22").Value = 5900
$ws.Range("L22").Value = 5900
$ws.Range("N22").Value = -6490

$ws.Range("H23").Value = 17500
$ws.Range("I23").Value = 17500
$ws.Range("K23").Value = 17500
$ws.Range("M23").Value = -17270

$ws.Range("H27").Value = 3010.75
$ws.Range("J27").Value = 5900
$ws.Range("L27").Value = 5900
$ws.Range("N27").Value = -6114

$ws.Range("H40").Value = 6016.533
$ws.Range("I40").Value = 5675
$ws.Range("K40").Value = 5675
$ws.Range("M40").Value = -5539

$ws.Range("H46").Value = 8747.333000000001
$ws.Range("I46").Value = 1492.5
$ws.Range("J46").Value = 12374.75
$ws.Range("K46").Value = 1492.5
$ws.Range("L46").Value = 12374.75
$ws.Range("M46").Value = -1304.5
$ws.Range("N46").Value = -12750.75

$ws.Range("H100").Value = 3531.375
$ws.Range("J100").Value = 3599
$ws.Range("L100").Value = 3599
$ws.Range("N100").Value = -4681

$ws.Range("H132").Value = 2277384.8
$ws.Range("I132").Value = 2461172.8
$ws.Range("J132").Value = 10666.333
$ws.Range("K132").Value = 7383518.399999999
$ws.Range("L132").Value = 31998.999
$ws.Range("M132").Value = -7380988.399999999
$ws.Range("N132").Value = -37058.999

$ws.Range("H136").Value = 9807999
$ws.Range("I136").Value = 15876156
$ws.Range("J136").Value = 5590.4614
$ws.Range("K136").Value = 47628468
$ws.Range("L136").Value = 16771.3842
$ws.Range("M136").Value = -47625918
$ws.Range("N136").Value = -21871.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents() | Out-Null
$ws.Range("N38").ClearContents() | Out-Null

$ws.Range("H44").Value = 37499.75
$ws.Range("I44").Value = 35999
$ws.Range("J44").Value = 38000
$ws.Range("K44").Value = 35999
$ws.Range("L44").Value = 38000
$ws.Range("M44").Value = -35445
$ws.Range("N44").Value = -39108

$ws.Range("H48").Value = 36249.75
$ws.Range("I48").Value = 36999
$ws.Range("J48").Value = 36000
$ws.Range("K48").Value = 36999
$ws.Range("L48").Value = 36000
$ws.Range("M48").Value = -36430
$ws.Range("N48").Value = -37138

$ws.Range("H126").Value = 102382790
$ws.Range("J126").Value = 1904
$ws.Range("L126").Value = 5712
$ws.Range("N126").Value = -10652

$ws.Range("H132").Value = 3371.1143
$ws.Range("I132").Value = 3378.4546
$ws.Range("K132").Value = 10135.3638
$ws.Range("M132").Value = -7605.363799999999

